$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template siswa: sisakan hanya satu baris contoh (baris 2) dan hapus
# baris-baris contoh tambahan (baris 3-6) sehingga template lebih ringkas.
$ws.Range("A3:E6").EntireRow.Delete()

# Jadikan baris contoh yang tersisa sebagai placeholder generik.
# Kolom NIS diformat sebagai Teks agar nilai numerik "1234" tetap
# tersimpan sebagai teks (konsisten dengan data NIS lain di kolom ini).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1234"
$ws.Range("C2").Value = "Nama Siswa"
